$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits before "Probably not,"
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) After "... as URIs." add two spaces, a new (fresh) "_GoBack"
#    bookmark, and a new sentence about CI.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("as URIs.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter("  ")
    $rng.Collapse(0)
    $rng.InsertAfter("That is, with some care, you can apply CI (Continuous Integration) even to APIs used by thousands of external parties!")

    # Re-seat a fresh, zero-length range at the boundary between the two
    # spaces and the new sentence, and drop a "_GoBack" bookmark there
    # (mirrors the boundary-style placement Word itself produces).
    $boundary = $rng.Start
    $bmRng = $d.Range($boundary, $boundary)
    $d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null
}

# ------------------------------------------------------------------
# 3) Split "be the same entity as the one " into two runs at
#    "be the same en" | "tity as the one ", and change the trailing
#    "!" to "." a couple of sentences later.
# ------------------------------------------------------------------
$rng2 = $d.Content
$foundSplit = $rng2.Find.Execute("be the same en", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundSplit) {
    $splitPos = $rng2.End
    $tmpBmRng = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("ZZTMP_SPLIT", $tmpBmRng) | Out-Null
    $d.Bookmarks("ZZTMP_SPLIT").Delete()
}

$rng3 = $d.Content
$foundBang = $rng3.Find.Execute("being published!", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundBang) {
    $bangRng = $d.Range($rng3.End - 1, $rng3.End)
    $bangRng.Text = "."
}

# ------------------------------------------------------------------
# 4) Adjust page margins: top 1134 -> 993 twips, bottom 993 -> 142 twips
#    (PageSetup works in points; 20 twips == 1 point).
# ------------------------------------------------------------------
$ps = $d.PageSetup
$ps.TopMargin = 993 / 20
$ps.BottomMargin = 142 / 20

Write-Output "edits applied"
